$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 32 appended to the results table.
# A32 ("2025-03-25") must stay literal text like the other Date cells, not be
# auto-parsed into a date serial, so format it as Text before assigning the
# value, then restore the "Normal" style so no stray number-format sticks to
# the cell.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2025-03-25"
$ws.Range("A32").Style = "Normal"

$ws.Range("B32").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C32").Value = "NA"
$ws.Range("D32").Value = 1
